$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select the header row, then delete it so the data shifts up by one row,
# matching the target layout (Pavan now on row 1 ... Urmila on row 6).
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()
